# Disponibilidad.xlsx - "Actualizar 02-06-2021 07-42-31" automated update
#
# 1) The 14 rows belonging to the previous check run (rows 1080-1093) get
#    their "Fecha" timestamp nudged forward a few seconds
#    (44233.2999336663 -> 44233.29993366898).
# 2) A brand-new check run is appended as 14 more rows (1094-1107), mirroring
#    the same Nombre/URL/Disponibilidad layout with a fresh timestamp
#    (44233.32113771647) and live hyperlinks in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) refresh the timestamp on the previous run's 14 rows -------------
for ($r = 1080; $r -le 1093; $r++) {
    $ws.Range("D$r").Value = 44233.29993366898
}

# --- 2) append the new run's 14 rows -------------------------------------
$items = @(
    @{Name="Odoo";              CellText="https://www.dataintelligence-group.com/";                       Address="https://www.dataintelligence-group.com/";                       SubAddress=""},
    @{Name="Blackbox";          CellText="https://serviciodashboard.azurewebsites.net/";                   Address="https://serviciodashboard.azurewebsites.net/";                   SubAddress=""},
    @{Name="PowerBI";           CellText="https://powerbi.microsoft.com/es-es/";                           Address="https://powerbi.microsoft.com/es-es/";                           SubAddress=""},
    @{Name="Dropbox";           CellText="https://www.dropbox.com/";                                       Address="https://www.dropbox.com/";                                       SubAddress=""},
    @{Name="Odoo";              CellText="https://dataintelligence.store/";                                Address="https://dataintelligence.store/";                                SubAddress=""},
    @{Name="GEE";                CellText="https://app-data-i.users.earthengine.app/";                      Address="https://app-data-i.users.earthengine.app/";                      SubAddress=""},
    @{Name="UtilidadesOdoo";    CellText="https://odooutil.azurewebsites.net/";                            Address="https://odooutil.azurewebsites.net/";                            SubAddress=""},
    @{Name="Filtros Dashboard"; CellText="https://filtradordashboard.azurewebsites.net/";                  Address="https://filtradordashboard.azurewebsites.net/";                  SubAddress=""},
    @{Name="MapStore";          CellText="https://ide.dataintelligence-group.com/mapstore/#/";             Address="https://ide.dataintelligence-group.com/mapstore/";               SubAddress="/"},
    @{Name="GeoServer";         CellText="https://ide.dataintelligence-group.com/geoserver/web/?0";        Address="https://ide.dataintelligence-group.com/geoserver/web/?0";        SubAddress=""},
    @{Name="Tomcat";            CellText="https://ide.dataintelligence-group.com/";                        Address="https://ide.dataintelligence-group.com/";                        SubAddress=""},
    @{Name="Shiny";             CellText="https://rpubs.com/dataintelligence/";                            Address="https://rpubs.com/dataintelligence/";                            SubAddress=""},
    @{Name="Github";            CellText="https://github.com/Sud-Austral/";                                Address="https://github.com/Sud-Austral/";                                SubAddress=""},
    @{Name="EZ Exporter";       CellText="https://ezexporter.highviewapps.com/exports/export-profile/";    Address="https://ezexporter.highviewapps.com/exports/export-profile/";    SubAddress=""}
)

$newTimestamp = 44233.32113771647
$startRow = 1094
$row = $startRow

foreach ($item in $items) {
    $a = $ws.Range("A$row")
    $b = $ws.Range("B$row")
    $c = $ws.Range("C$row")
    $d = $ws.Range("D$row")

    $a.Value = $item.Name
    $b.Value = $item.CellText

    if ($item.SubAddress -ne "") {
        $ws.Hyperlinks.Add($b, $item.Address, $item.SubAddress)
    } else {
        $ws.Hyperlinks.Add($b, $item.Address)
    }
    $b.Style = "Hyperlink"

    $c.Value = "Disponible"

    $d.Value = $newTimestamp
    $d.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $row = $row + 1
}

Write-Output "Updated rows 1080-1093 and appended rows $startRow-$($row - 1)."
